$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.636.38"
$ws.Range("E2").Value = "  -2.25%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.760.46"
$ws.Range("E3").Value = "  -2.98%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.67"
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4354"
$ws.Range("E7").Value = "  +0.34%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3608"
$ws.Range("E8").Value = "  -1.34%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07571"
$ws.Range("E9").Value = "  -1.38%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.06"
$ws.Range("E10").Value = "  -6.24%  "
$ws.Range("E11").Value = "  -2.60%  "
$ws.Range("E12").Value = "  -0.03%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.84"
$ws.Range("E13").Value = "  -5.22%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.079"
$ws.Range("E14").Value = "  -3.30%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.239"
$ws.Range("E15").Value = "  -3.65%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.757.02"
$ws.Range("E16").Value = "  -4.01%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "92.70"
$ws.Range("E17").Value = "  -0.94%  "
$ws.Range("E18").Value = "  -1.08%  "
$ws.Range("E19").Value = "  -2.20%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.0000"
$ws.Range("E20").Value = "  -0.11%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.11"
$ws.Range("E21").Value = "  -2.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.851"
$ws.Range("E22").Value = "  -6.43%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.666.71"
$ws.Range("E23").Value = "  -2.23%  "
$ws.Range("E24").Value = "  -2.45%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.093"
$ws.Range("E25").Value = "  +1.73%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "162.82"
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("E27").Value = "  -0.77%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.957.35"
$ws.Range("E28").Value = "  -3.86%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.145"
$ws.Range("E29").Value = "  -6.18%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.07"
$ws.Range("E30").Value = "  -2.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.101"
$ws.Range("E31").Value = "  -10.05%  "
$ws.Range("B32").Value = "HuobiToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.664"
$ws.Range("E32").Value = "  +5.49%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.587"
$ws.Range("E33").Value = "  -6.07%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08983"
$ws.Range("E34").Value = "  -1.94%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "12.20"
$ws.Range("E35").Value = "  -5.92%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6443"
$ws.Range("E37").Value = "  -1.74%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2109"
$ws.Range("E38").Value = "  -2.94%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06020"
$ws.Range("E39").Value = "  -2.65%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.957"
$ws.Range("E40").Value = "  -4.48%  "
$ws.Range("E41").Value = "  -0.60%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.001"
$ws.Range("E42").Value = "  -0.04%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.401"
$ws.Range("E43").Value = "  -1.58%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.899"
$ws.Range("E44").Value = "  -2.52%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.32"
$ws.Range("E45").Value = "  -4.00%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5943"
$ws.Range("E46").Value = "  -2.57%  "
$ws.Range("E47").Value = "  -0.88%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.992"
$ws.Range("E48").Value = "  -1.24%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "122.41"
$ws.Range("E49").Value = "  -2.52%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.170"
$ws.Range("E50").Value = "  +1.26%  "
$ws.Range("E51").Value = "  -1.76%  "
